$wb = $excel.ActiveWorkbook
$tmp = $wb.Worksheets.Add()
